$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-9 down to 5-10.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record's data.
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44603
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100112001
$ws.Cells.Item(4, 7).Value = "Berenjena"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 140
$ws.Cells.Item(4, 11).Value = 5500
$ws.Cells.Item(4, 12).Value = 6000
$ws.Cells.Item(4, 13).Value = 5750
$ws.Cells.Item(4, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 96
$ws.Cells.Item(4, 17).Value = 60
$ws.Cells.Item(4, 18).Value = "Hortaliza"

$wb.Save()
